$wb = $excel.ActiveWorkbook

$wsBasic = $wb.Worksheets.Item("Basic Game rubric")

# --- Content fixes (pixl inventory visual bug) ---

# "Animations" row: score 2 -> 3
$wsBasic.Range("B3").Value = 3

# "Game implementation" row: score 2 -> 3, and update description text
$wsBasic.Range("B5").Value = 3
$wsBasic.Range("C5").Value = "Player can move around, enemies walk around and attack, player can change attack, use items, advance to next level, interact with objects"

# --- Selection / active sheet state ---
# The previously active sheet ("Game extras") is no longer selected;
# "Basic Game rubric" becomes the active sheet with C6 selected.
$wsBasic.Activate()
$wsBasic.Range("C6").Select()
